$wb = $excel.ActiveWorkbook
$wsLista = $wb.Worksheets.Item("Lista de Tabelas")
$wsC4 = $wb.Worksheets.Item("C4 - Massa de rendimentos")

# --- Update the table title text (rolling the reference year from 2014 to 2015) ---
$novoTitulo = "Tabela 1 - Evolução do Índice de Concentração C4 do Total de Rendimentos dos das Atividades Relacionadas por UF entre 2007 e 2015"

$wsLista.Range("A2").Value = $novoTitulo
$wsC4.Range("A1").Value = $novoTitulo

# --- Shift the year headers on the data sheet one year forward (2010-2014 -> 2011-2015) ---
$wsC4.Range("E2").Value = 2011
$wsC4.Range("F2").Value = 2012
$wsC4.Range("G2").Value = 2013
$wsC4.Range("H2").Value = 2014
$wsC4.Range("I2").Value = 2015

# --- Update the underlying C4 concentration values for the shifted years ---
$wsC4.Range("B3").Value = 0.64941377886383655
$wsC4.Range("C3").Value = 0.66888748756939687
$wsC4.Range("D3").Value = 0.63269189492511901
$wsC4.Range("E3").Value = 0.64001049055683878
$wsC4.Range("F3").Value = 0.62574470603802668
$wsC4.Range("G3").Value = 0.64228437110822345
$wsC4.Range("H3").Value = 0.64233690451414227

# --- Update view/selection state: the C4 sheet becomes the active tab ---
$wsLista.Range("H10").Select()
$wsC4.Activate()
$wsC4.Range("J8").Select()
